$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Force the cell to stay a text string even when the literal
    # looks like a number (e.g. "526.24"), then drop back to the
    # workbook default style so no stray formatting is introduced.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# --- Row-by-row price (D) and volume/% (E) updates ---
Set-TextValue $ws.Range("D2") "57.859.09"
Set-TextValue $ws.Range("E2") "  +0.57%  "

Set-TextValue $ws.Range("D3") "3.130.79"
Set-TextValue $ws.Range("E3") "  +2.04%  "

Set-TextValue $ws.Range("E4") "  +0.00%  "

Set-TextValue $ws.Range("D5") "526.24"
Set-TextValue $ws.Range("E5") "  +1.58%  "

Set-TextValue $ws.Range("D6") "141.78"
Set-TextValue $ws.Range("E6") "  +0.53%  "

Set-TextValue $ws.Range("D8") "3.130.53"
Set-TextValue $ws.Range("E8") "  +1.99%  "

Set-TextValue $ws.Range("E9") "  -0.15%  "

Set-TextValue $ws.Range("D10") "7.28"
Set-TextValue $ws.Range("E10") "  +0.65%  "

Set-TextValue $ws.Range("E11") "  +2.12%  "

Set-TextValue $ws.Range("D12") "0.386"
Set-TextValue $ws.Range("E12") "  +3.03%  "

Set-TextValue $ws.Range("D13") "3.672.78"
Set-TextValue $ws.Range("E13") "  +2.07%  "

Set-TextValue $ws.Range("D14") "0.132"
Set-TextValue $ws.Range("E14") "  +1.73%  "

Set-TextValue $ws.Range("D15") "26.38"
Set-TextValue $ws.Range("E15") "  +3.12%  "

Set-TextValue $ws.Range("E16") "  +1.46%  "

Set-TextValue $ws.Range("D17") "57.958.40"
Set-TextValue $ws.Range("E17") "  +0.59%  "

Set-TextValue $ws.Range("D18") "3.139.32"
Set-TextValue $ws.Range("E18") "  +2.34%  "

Set-TextValue $ws.Range("E19") "  +1.01%  "

Set-TextValue $ws.Range("D20") "12.91"
Set-TextValue $ws.Range("E20") "  +0.38%  "

Set-TextValue $ws.Range("D21") "8.11"
Set-TextValue $ws.Range("E21") "  +0.10%  "

Set-TextValue $ws.Range("D22") "337.98"
Set-TextValue $ws.Range("E22") "  +1.81%  "

Set-TextValue $ws.Range("D23") "0.998"
Set-TextValue $ws.Range("E23") "  -0.13%  "

Set-TextValue $ws.Range("D24") "0.513"
Set-TextValue $ws.Range("E24") "  +2.86%  "

Set-TextValue $ws.Range("D25") "66.91"
Set-TextValue $ws.Range("E25") "  +1.83%  "

Set-TextValue $ws.Range("D26") "0.169"
Set-TextValue $ws.Range("E26") "  -0.01%  "

Set-TextValue $ws.Range("E27") "  -0.02%  "

Set-TextValue $ws.Range("D28") "0.0₃0933"
Set-TextValue $ws.Range("E28") "  +3.56%  "

Set-TextValue $ws.Range("E29") "  +4.15%  "

Set-TextValue $ws.Range("D30") "0.998"
Set-TextValue $ws.Range("E30") "  +0.00%  "

Set-TextValue $ws.Range("D31") "7.26"
Set-TextValue $ws.Range("E31") "  +0.85%  "

Set-TextValue $ws.Range("D32") "1.88"
Set-TextValue $ws.Range("E32") "  +3.38%  "

Set-TextValue $ws.Range("E33") "  +2.25%  "

Set-TextValue $ws.Range("D34") "21.03"
Set-TextValue $ws.Range("E34") "  +1.34%  "

Set-TextValue $ws.Range("D35") "155.78"
Set-TextValue $ws.Range("E35") "  +0.67%  "

Set-TextValue $ws.Range("D36") "4.69"
Set-TextValue $ws.Range("E36") "  +4.30%  "

Set-TextValue $ws.Range("D37") "6.13"
Set-TextValue $ws.Range("E37") "  +3.67%  "

Set-TextValue $ws.Range("D38") "27.17"
Set-TextValue $ws.Range("E38") "  +0.44%  "

Set-TextValue $ws.Range("D39") "1.30"
Set-TextValue $ws.Range("E39") "  +2.88%  "

Set-TextValue $ws.Range("E40") "  -0.50%  "

Set-TextValue $ws.Range("D41") "3.175.57"

# Rows 42/43 swap places (Mantle <-> Stacks) with refreshed figures
Set-TextValue $ws.Range("B42") "Stacks"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D42") "1.54"
Set-TextValue $ws.Range("E42") "  +13.03%  "

Set-TextValue $ws.Range("B43") "Mantle"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D43") "0.691"
Set-TextValue $ws.Range("E43") "  +5.76%  "

Set-TextValue $ws.Range("D44") "3.93"
Set-TextValue $ws.Range("E44") "  +0.27%  "

Set-TextValue $ws.Range("D45") "37.00"
Set-TextValue $ws.Range("E45") "  +0.64%  "

Set-TextValue $ws.Range("E46") "  -0.03%  "

Set-TextValue $ws.Range("D47") "2.302.68"
Set-TextValue $ws.Range("E47") "  +1.78%  "

Set-TextValue $ws.Range("D48") "0.0261"
Set-TextValue $ws.Range("E48") "  +0.96%  "

Set-TextValue $ws.Range("E49") "  +7.84%  "

Set-TextValue $ws.Range("D50") "21.12"
Set-TextValue $ws.Range("E50") "  +1.13%  "

Set-TextValue $ws.Range("D51") "6.02"
Set-TextValue $ws.Range("E51") "  +2.53%  "
